$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 19 (pushes old row19 -> row20, old rows24/25 -> 25/26),
#    then clone the formatting of row 18 (regular data-row style) into the new row 19.
$ws.Rows("19").Insert()
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2. Update the "VALOR MORA" total and "Cant. Periodos" count.
$ws.Range("E11").Value = 279760
$ws.Range("F13").Value = 5

# 3. Update the period/value columns for the 5 worker-period rows.
#    Row 16
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500
#    Row 17
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500
#    Row 18
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500
#    Row 19 (newly inserted row, same style as 16-18)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45508682"
$ws.Range("D19").Value = "ANGELA MARIA VALDEZ BOHORQUEZ"
$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500
#    Row 20 (previously row 19, carries the "last row" bottom-border style)
$ws.Range("E20").Value = "2503"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1423500

Write-Host "Edit applied"
